# WorkLog-PiersenSchuiling.xlsx update
# - Fill in the timecard entries for rows 68-73 (new work log rows)
# - Update the active selection to N66
# (row-height / default-col-width cosmetics are left to the host app's
#  own re-layout, since they are not explicit user edits)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of timecard data to add. Row 60 (A60:E60) already carries the
# exact date / time / duration / description number formats (styles
# 13, 17, 17, 2, 3) that these new rows should use, so we copy that
# formatting across before writing the new values.
$newRows = @(
    @{ row = 68; date = 45769; start = 0.875;             finish = 0.916666666666667;  dur = 0.0416666666666667; desc = "Finding seg fault in Test::followPath" },
    @{ row = 69; date = 45772; start = 0.625;              finish = 0.708333333333333;  dur = 0.0833333333333333; desc = "Configuring project for Windows" },
    @{ row = 70; date = 45772; start = 0.791666666666667;  finish = 0.822916666666667;  dur = 0.03125;            desc = "Finding seg fault in Test::followPath" },
    @{ row = 71; date = 45772; start = 0.875;              finish = 0.9375;             dur = 0.0625;             desc = "Correcting incident edges to detect proper edge and road segment in Test::followPath" },
    @{ row = 72; date = 45775; start = 0.75;                finish = 0.833333333333333;  dur = 0.0833333333333333; desc = "Finishing final tests and filling out Test Results Document" },
    @{ row = 73; date = 45775; start = 0.875;              finish = 0.979166666666667;  dur = 0.104166666666667;  desc = "Finishing Testing Document, User Documentation, and remaining final submissions" }
)

$formatDonor = $ws.Range("A60:E60")

foreach ($entry in $newRows) {
    $r = $entry.row

    $formatDonor.Copy()
    $ws.Range("A" + $r + ":E" + $r).PasteSpecial(-4122)

    $ws.Range("A" + $r).Value = $entry.date
    $ws.Range("B" + $r).Value = $entry.start
    $ws.Range("C" + $r).Value = $entry.finish
    $ws.Range("D" + $r).Value = $entry.dur
    $ws.Range("E" + $r).Value = $entry.desc
}

# Move the active selection / view to reflect where the author ended up
$ws.Range("N66").Select()

Write-Host "Updated rows 68-73 and selection"
